# Update codes with alert: replace the second contact's data, add three
# more contacts with hyperlinked emails, and pad the sheet with repeated
# name/code rows down to row 15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update existing row 2 (name/code change; e-mail + hyperlink change) ---
$ws.Range("A2").Value = "arnolf"
$ws.Range("B2").Value = "jb"

# Row 2's hyperlink needs to point at the new address. Remove the old
# hyperlink object tied to C2 (only it, not the sibling on C1) before
# writing the new address + re-adding the link.
foreach ($h in @($ws.Hyperlinks)) {
    if ($h.Range.Row -eq 2) {
        $h.Delete()
    }
}
$ws.Range("C2").Value = "g2@mailinator.com"
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:g2@mailinator.com", [Type]::Missing, [Type]::Missing, "g2@mailinator.com")
# Re-adding a hyperlink re-applies Excel's built-in "Hyperlink" cell style;
# paste the plain formatting from A1 back over it so it keeps the workbook's
# original (unstyled) look.
$ws.Range("A1").Copy()
$ws.Range("C2").PasteSpecial(-4122)

# --- New rows 3 & 4: same name/code, new hyperlinked e-mail addresses ---
$ws.Range("A3").Value = "arnolf"
$ws.Range("B3").Value = "jb"
$ws.Range("C3").Value = "gm@mailinator.com"
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:gm@mailinator.com", [Type]::Missing, [Type]::Missing, "gm@mailinator.com")
$ws.Range("A1").Copy()
$ws.Range("C3").PasteSpecial(-4122)

$ws.Range("A4").Value = "arnolf"
$ws.Range("B4").Value = "jb"
$ws.Range("C4").Value = "gh@mailinator.com"
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:gh@mailinator.com", [Type]::Missing, [Type]::Missing, "gh@mailinator.com")
$ws.Range("A1").Copy()
$ws.Range("C4").PasteSpecial(-4122)

# --- Rows 5-15: just the repeated name/code pair, no e-mail ---
for ($r = 5; $r -le 15; $r++) {
    $ws.Cells.Item($r, 1).Value = "arnolf"
    $ws.Cells.Item($r, 2).Value = "jb"
}

[void]$ws.Range("D4").Select()
